$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new weekly data rows above the current row 59, pushing the
# existing data (rows 59-82) down to rows 61-84.
$ws.Rows.Item(59).Insert()
$ws.Rows.Item(59).Insert()

# New row 59: Cilantro "Primera" reading for 2022-09-05
$ws.Cells.Item(59, 1).Value2 = 7
$ws.Cells.Item(59, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(59, 3).Value2 = "Ñuble"
$ws.Cells.Item(59, 4).Value2 = 44809
$ws.Cells.Item(59, 5).Value2 = 16
$ws.Cells.Item(59, 6).Value2 = 100112040
$ws.Cells.Item(59, 7).Value2 = "Cilantro"
$ws.Cells.Item(59, 8).Value2 = "Sin especificar"
$ws.Cells.Item(59, 9).Value2 = "Primera"
$ws.Cells.Item(59, 10).Value2 = 240
$ws.Cells.Item(59, 11).Value2 = 700
$ws.Cells.Item(59, 12).Value2 = 800
$ws.Cells.Item(59, 13).Value2 = 750
$ws.Cells.Item(59, 14).Value2 = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(59, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(59, 16).Value2 = 750
$ws.Cells.Item(59, 17).Value2 = 1
$ws.Cells.Item(59, 18).Value2 = "Hortaliza"

# New row 60: Cilantro "Segunda" reading for 2022-09-05
$ws.Cells.Item(60, 1).Value2 = 7
$ws.Cells.Item(60, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(60, 3).Value2 = "Ñuble"
$ws.Cells.Item(60, 4).Value2 = 44809
$ws.Cells.Item(60, 5).Value2 = 16
$ws.Cells.Item(60, 6).Value2 = 100112040
$ws.Cells.Item(60, 7).Value2 = "Cilantro"
$ws.Cells.Item(60, 8).Value2 = "Sin especificar"
$ws.Cells.Item(60, 9).Value2 = "Segunda"
$ws.Cells.Item(60, 10).Value2 = 150
$ws.Cells.Item(60, 11).Value2 = 600
$ws.Cells.Item(60, 12).Value2 = 600
$ws.Cells.Item(60, 13).Value2 = 600
$ws.Cells.Item(60, 14).Value2 = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(60, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(60, 16).Value2 = 600
$ws.Cells.Item(60, 17).Value2 = 1
$ws.Cells.Item(60, 18).Value2 = "Hortaliza"

$ws.Cells.Item(59, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(60, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
